$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.293.91'
$ws.Range("E2").Value = '  +1.24%  '

$ws.Range("D3").Value = '2.833.02'
$ws.Range("E3").Value = '  +3.27%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''357.43'
$ws.Range("E5").Value = '  +7.53%  '

$ws.Range("D6").Value = '''114.56'
$ws.Range("E6").Value = '  -1.54%  '

$ws.Range("D7").Value = '''0.549'
$ws.Range("E7").Value = '  +2.99%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = '''0.608'
$ws.Range("E9").Value = '  +6.26%  '

$ws.Range("D10").Value = '''42.08'
$ws.Range("E10").Value = '  +1.43%  '

$ws.Range("D11").Value = '''0.0851'
$ws.Range("E11").Value = '  +2.01%  '

$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("E13").Value = '  +1.51%  '

$ws.Range("D14").Value = '''7.79'
$ws.Range("E14").Value = '  +3.22%  '

$ws.Range("D15").Value = '3.285.65'
$ws.Range("E15").Value = '  +3.47%  '

$ws.Range("D16").Value = '2.827.63'
$ws.Range("E16").Value = '  +2.81%  '

$ws.Range("E17").Value = '  +1.63%  '

$ws.Range("D18").Value = '52.239.86'
$ws.Range("E18").Value = '  +1.25%  '

$ws.Range("D19").Value = '''3.17'
$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("D20").Value = '''7.29'
$ws.Range("E20").Value = '  +7.15%  '

$ws.Range("D21").Value = '''13.75'
$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("E22").Value = '  +3.05%  '

$ws.Range("D23").Value = '''271.77'
$ws.Range("E23").Value = '  -2.82%  '

$ws.Range("D24").Value = '''69.80'
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("E25").Value = '  +6.38%  '

$ws.Range("D26").Value = '''26.81'
$ws.Range("E26").Value = '  +0.43%  '

$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("E28").Value = '  +0.87%  '

$ws.Range("E29").Value = '  +1.30%  '

$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("D31").Value = '''50.65'
$ws.Range("E31").Value = '  +0.99%  '

$ws.Range("D32").Value = '''33.76'
$ws.Range("E32").Value = '  -3.29%  '

$ws.Range("D33").Value = '''5.91'
$ws.Range("E33").Value = '  +6.68%  '

$ws.Range("D34").Value = '''0.0441'
$ws.Range("E34").Value = '  +27.89%  '

$ws.Range("D35").Value = '''0.0832'
$ws.Range("E35").Value = '  +1.81%  '

$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("E37").Value = '  +1.32%  '

$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = '''18.52'
$ws.Range("E38").Value = '  -2.27%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''4.89'
$ws.Range("E39").Value = '  -0.61%  '

$ws.Range("D40").Value = '''3.20'
$ws.Range("E40").Value = '  +1.97%  '

$ws.Range("D41").Value = '''23.66'
$ws.Range("E41").Value = '  +3.13%  '

$ws.Range("D42").Value = '''2.57'
$ws.Range("E42").Value = '  +8.75%  '

$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''127.41'
$ws.Range("E43").Value = '  -1.16%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.115'
$ws.Range("E44").Value = '  +2.00%  '

$ws.Range("E45").Value = '  +2.30%  '

$ws.Range("D46").Value = '''3.37'
$ws.Range("E46").Value = '  +1.85%  '

$ws.Range("D47").Value = '2.047.30'
$ws.Range("E47").Value = '  -2.42%  '

$ws.Range("E48").Value = '  +3.78%  '

$ws.Range("D49").Value = '''0.972'
$ws.Range("E49").Value = '  +13.07%  '

$ws.Range("D50").Value = '''5.73'
$ws.Range("E50").Value = '  +3.76%  '

$ws.Range("D51").Value = '''60.32'
$ws.Range("E51").Value = '  +0.76%  '
